$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5569
$ws1.Range("F8").Value = 915
$ws1.Range("F10").Value = 2493
$ws1.Range("F11").Value = 83
$ws1.Range("F12").Value = 113
$ws1.Range("F14").Value = 73
$ws1.Range("F16").Value = 2330
$ws1.Range("F17").Value = 285

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 101

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5569
$ws4.Range("F6").Value = 101
$ws4.Range("F10").Value = 915
$ws4.Range("F12").Value = 2493
$ws4.Range("F13").Value = 83
$ws4.Range("F14").Value = 113
$ws4.Range("F17").Value = 73
$ws4.Range("F19").Value = 2330
$ws4.Range("F20").Value = 285
